# Add a new row of data ("Pesho", 100) below the existing table data
# (row 5 -> row 6), matching the commit's "Some fixes in DB homeworks"
# change to myWorksheet.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new values.
$ws.Range("A6").Value = "Pesho"
$ws.Range("B6").Value = 100

# The previous last row (B5) carries the "bottom-right corner" cell
# border style; copy that same formatting onto the new A6:B6 cells so
# they pick up style index 6 for both columns, exactly like the target
# workbook (the Excel table itself is left at its original A1:B5 range).
$ws.Range("B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4122)
